# MPSTOOLS-1 - Finished Quotegen Models
# Collapse the "Hardware Device List" sheet down to a single device block
# (Device 1) and replace its placeholder item rows with the finished
# HP LaserJet M9050 line items.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Device 2" / "Device 3" / "Device 4" blocks (rows 8-23), then
# remove the now-trailing blank separator row (row 7) so that what used to
# be row 24 (the second item row of the old "Device 1" block) shifts up to
# become the new row 7.
$ws.Rows("8:23").Delete()
$ws.Rows("7:7").Delete()

# Replace the placeholder device/items with the finished model data.
$ws.Range("A5").Value = "HP LaserJet M9050"
$ws.Range("B5").Value = "CC395A"
$ws.Range("A6").Value = "     8-Bin Mailbox"
$ws.Range("B6").Value = "Q5693A"
$ws.Range("A7").Value = "HP LaserJet MFP 3000-sheet Stapler/Stacker"
$ws.Range("B7").Value = "C8085A"

# Restore the selection to the first item row.
$ws.Range("A5").Select()
